# Regenerate save_data to use K (strikeouts) instead of Strike# (old "Strike#"
# derived value). This recalculates column G ("K") for each outing row and
# writes the new strikeout totals (s_vals) in place of the previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2..37, replacing the old Strike# figures.
$newK = @{
    2  = 4
    3  = 6
    4  = 10
    5  = 5
    6  = 7
    7  = 8
    8  = 6
    9  = 7
    10 = 8
    11 = 2
    12 = 7
    13 = 7
    14 = 9
    15 = 8
    16 = 6
    17 = 4
    18 = 5
    19 = 6
    20 = 7
    21 = 7
    22 = 7
    23 = 5
    24 = 3
    25 = 5
    26 = 11
    27 = 3
    28 = 2
    29 = 5
    30 = 3
    31 = 4
    32 = 7
    33 = 5
    34 = 0
    35 = 5
    36 = 2
    37 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
